# "Users and Roles chapter compleated"
# Mark the "users-and-roles.md" row (row 35) as Completed / Converted / Checked,
# matching the green highlight style already used for other finished rows
# (e.g. row 5, row 7, row 17 ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list")

# Row 35 = users-and-roles.md -> B (Completed), C (Converted), D (Checked)
$rng = $ws.Range("B35:D35")
$rng.Value = "x"

# Apply the same "done" look used elsewhere in the sheet: centered text on a
# green fill (style index 2 in styles.xml - fill 92D050 + center alignment).
$rng.HorizontalAlignment = -4108   # xlCenter
$rng.Interior.Color = 5296274      # RGB(0x92,0xD0,0x50) in BGR OLE order

# Reset the sheet's saved cursor/selection back to the top-left cell.
$ws.Range("A1").Select()

Write-Host "Marked users-and-roles.md as completed/converted/checked"
